$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '58.821.87'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.155.70'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '532.30'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '140.30'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.96%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.537'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +16.30%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.34'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.61%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +5.80%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +3.41%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.71%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.700.21'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.22'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +6.02%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '58.835.98'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.27'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +4.40%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.144.65'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.05'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.19'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '373.05'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +5.37%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.80'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.16'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.520'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.14%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.168'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.29'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +13.74%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0865'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.43%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.16'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +3.83%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.13'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.36%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +3.85%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.80%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '159.01'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.29'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +7.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '25.23'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.03%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0686'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.78%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.636.07'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +9.86%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +5.75%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +7.78%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '38.82'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +3.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.713'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.196.72'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +14.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.987'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.74%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '20.24'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.56%  '
